# Model fixes to elec sector BECF, BPMCCS, BCRbQ
# Updates underlying input data on "water & waste" and "BNRbI" sheets.
# All dependent formulas recalc automatically.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "water & waste": update 2017 state population (E5) and the
# U.S. population projection series (C17:AJ17).
# ---------------------------------------------------------------------
$wsWW = $wb.Worksheets.Item("water & waste")

$wsWW.Range("E5").Value = 3012.223

$popVals = @(3012.223,3056.748,3101.2719999999999,3145.797,3190.3220000000001,3234.846,3279.3710000000001,3323.895,3368.42,3412.9450000000002,3457.4690000000001,3501.9940000000001,3546.518,3591.0430000000001,3635.5680000000002,3680.0920000000001,3724.6170000000002,3769.1410000000001,3813.6660000000002,3858.1909999999998,3902.7150000000001,3947.24,3991.7640000000001,4036.2890000000002,4080.8139999999999,4125.3379999999997,4169.8630000000003,4214.3869999999997,4258.9120000000003,4303.4369999999999,4347.9610000000002,4392.4859999999999,4437.01,4481.5349999999999)

$popArr = New-Object 'object[,]' 1,$popVals.Length
for ($i = 0; $i -lt $popVals.Length; $i++) {
    $popArr[0,$i] = $popVals[$i]
}
$wsWW.Range("C17:AJ17").Value = $popArr

# ---------------------------------------------------------------------
# Sheet "BNRbI": update hardcoded nonfuel-revenue-by-industry rows.
# ---------------------------------------------------------------------
$wsBN = $wb.Worksheets.Item("BNRbI")

$cementVals = @(160546239.53999999,158560846.11000001,159571038.96000001,163296215.72999999,166672543.59,169641262.97999999,172619092.97999999,175104440.81999999,177597117,180272348.55000001,182953978.56,185666349.96000001,188444621.06999999,191212895.97,194037937.83000001,197215304.58000001,199902126.41999999,202330452.69,204829073.46000001,207066774.33000001,208823074.11000001,211685831.46000001,214930625.58000001,217521105.21000001,220643011.25999999,224268978.69,228067317.09,231660030.24000001,235114191.27000001,238637539.80000001,242016646.22999999,245413475.72999999)

$steelVals = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

$chemVals = @(795310233.38999999,811509478.17999995,843788549.38999999,871283171.58000004,899577401.09000003,921530389.72000003,938372979.10000014,952604303.19999993,966291655.33000004,981440904.63999987,996238516.19999993,1013149268.37,1027483552.51,1041094014.03,1054906753.91,1069330191.41,1084411610.97,1097494119.02,1110389596.1500001,1123598154.4100001,1137626096.05,1152084553.4300001,1165636023.4300001,1182417519.8800001,1196550522.97,1212250703.7,1228100024.8099999,1246191076.1700001,1261764433.8199999,1277371707.25,1295821687.04,1308990696.8)

$agVals = @(1129019924.6400001,1142852892.28,1162732484.5,1184366564.7,1201341265.5,1220003656.9000001,1238984452.5599999,1256818699.02,1274387231.8199999,1292199124.5599999,1310491086,1329930928.6199999,1351449113.48,1372356569.1600001,1392482451.24,1413769383.8800001,1434000248.22,1452826298.76,1470763530.5999999,1488571699.9000001,1506146086,1524807014.6199999,1543827726.0799999,1562499166.6600001,1581930890.96,1601844148.9200001,1622309722.96,1642900237.1600001,1663580567.3599999,1684476181.28,1705679421.54,1727114980.04)

$otherVals = @(10466673472.57,10435029779.02,10441886399.49,10598565788.02,10721640167,10864287609.77,10996864691.190001,11145417922.52,11332751415.93,11537818092.57,11751395908.42,11969289215.450001,12181999158,12405869999.950001,12647925109.75,12911787997.639999,13166743964.43,13406535939.51,13659925866.51,13900115965.26,14141371327.09,14403262231.540001,14677246686.9,14931654195.57,15206661734.879999,15480544316.790001,15745808281.42,15999134050.27,16256010231.51,16506376733.91,16773872035.620001,17049363901.99)

function Set-RowValues($ws, [string]$rowRange, $values) {
    $arr = New-Object 'object[,]' 1,$values.Length
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $ws.Range($rowRange).Value = $arr
}

Set-RowValues $wsBN "C2:AH2" $cementVals
Set-RowValues $wsBN "C4:AH4" $steelVals
Set-RowValues $wsBN "C5:AH5" $chemVals
Set-RowValues $wsBN "C8:AH8" $agVals
Set-RowValues $wsBN "C9:AH9" $otherVals

# ---------------------------------------------------------------------
# Active sheet / selection bookkeeping to match the saved workbook view.
# ---------------------------------------------------------------------
$wsBN.Select()
$wsBN.Range("C15").Select()
